$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = '62-38=24'
$t.Cell(1,2).Range.Text = '38-21=17'
$t.Cell(1,3).Range.Text = '0+35=35'
$t.Cell(1,4).Range.Text = '16-12=4'
$t.Cell(1,5).Range.Text = '50-5=45'
$t.Cell(2,1).Range.Text = '16+33=49'
$t.Cell(2,2).Range.Text = '45-21=24'
$t.Cell(2,3).Range.Text = '10+84=94'
$t.Cell(2,4).Range.Text = '0+64=64'
$t.Cell(2,5).Range.Text = '66-26=40'
$t.Cell(3,1).Range.Text = '35+2=37'
$t.Cell(3,2).Range.Text = '29+24=53'
$t.Cell(3,3).Range.Text = '32-12=20'
$t.Cell(3,4).Range.Text = '2+22=24'
$t.Cell(3,5).Range.Text = '71-7=64'
$t.Cell(4,1).Range.Text = '58+9=67'
$t.Cell(4,2).Range.Text = '11+88=99'
$t.Cell(4,3).Range.Text = '24+24=48'
$t.Cell(4,4).Range.Text = '19-11=8'
$t.Cell(4,5).Range.Text = '93-46=47'
$t.Cell(5,1).Range.Text = '88-3=85'
$t.Cell(5,2).Range.Text = '68-33=35'
$t.Cell(5,3).Range.Text = '13+69=82'
$t.Cell(5,4).Range.Text = '11+26=37'
$t.Cell(5,5).Range.Text = '91-38=53'
$t.Cell(6,1).Range.Text = '21-11=10'
$t.Cell(6,2).Range.Text = '83-60=23'
$t.Cell(6,3).Range.Text = '68-52=16'
$t.Cell(6,4).Range.Text = '1+45=46'
$t.Cell(6,5).Range.Text = '53+17=70'
$t.Cell(7,1).Range.Text = '83-71=12'
$t.Cell(7,2).Range.Text = '45+35=80'
$t.Cell(7,3).Range.Text = '22-11=11'
$t.Cell(7,4).Range.Text = '73-15=58'
$t.Cell(7,5).Range.Text = '75-18=57'
$t.Cell(8,1).Range.Text = '88-9=79'
$t.Cell(8,2).Range.Text = '2+22=24'
$t.Cell(8,3).Range.Text = '3+50=53'
$t.Cell(8,4).Range.Text = '9+84=93'
$t.Cell(8,5).Range.Text = '81+8=89'
$t.Cell(9,1).Range.Text = '7+22=29'
$t.Cell(9,2).Range.Text = '80+19=99'
$t.Cell(9,3).Range.Text = '16+29=45'
$t.Cell(9,4).Range.Text = '36-15=21'
$t.Cell(9,5).Range.Text = '86-54=32'
$t.Cell(10,1).Range.Text = '45+26=71'
$t.Cell(10,2).Range.Text = '4+93=97'
$t.Cell(10,3).Range.Text = '65-27=38'
$t.Cell(10,4).Range.Text = '25+57=82'
$t.Cell(10,5).Range.Text = '24+9=33'
$t.Cell(11,1).Range.Text = '88-3=85'
$t.Cell(11,2).Range.Text = '98-93=5'
$t.Cell(11,3).Range.Text = '62-34=28'
$t.Cell(11,4).Range.Text = '47-8=39'
$t.Cell(11,5).Range.Text = '3+65=68'
$t.Cell(12,1).Range.Text = '65+10=75'
$t.Cell(12,2).Range.Text = '41-19=22'
$t.Cell(12,3).Range.Text = '96+0=96'
$t.Cell(12,4).Range.Text = '91-19=72'
$t.Cell(12,5).Range.Text = '94-34=60'
$t.Cell(13,1).Range.Text = '87-13=74'
$t.Cell(13,2).Range.Text = '29+64=93'
$t.Cell(13,3).Range.Text = '70-15=55'
$t.Cell(13,4).Range.Text = '19-7=12'
$t.Cell(13,5).Range.Text = '29+59=88'
$t.Cell(14,1).Range.Text = '48+24=72'
$t.Cell(14,2).Range.Text = '35-30=5'
$t.Cell(14,3).Range.Text = '38+6=44'
$t.Cell(14,4).Range.Text = '55-52=3'
$t.Cell(14,5).Range.Text = '40+6=46'
$t.Cell(15,1).Range.Text = '49-8=41'
$t.Cell(15,2).Range.Text = '0+14=14'
$t.Cell(15,3).Range.Text = '67-48=19'
$t.Cell(15,4).Range.Text = '72-4=68'
$t.Cell(15,5).Range.Text = '45-36=9'
$t.Cell(16,1).Range.Text = '11+63=74'
$t.Cell(16,2).Range.Text = '4+36=40'
$t.Cell(16,3).Range.Text = '42-23=19'
$t.Cell(16,4).Range.Text = '57-51=6'
$t.Cell(16,5).Range.Text = '19+9=28'
$t.Cell(17,1).Range.Text = '27+22=49'
$t.Cell(17,2).Range.Text = '46+8=54'
$t.Cell(17,3).Range.Text = '39+41=80'
$t.Cell(17,4).Range.Text = '67-63=4'
$t.Cell(17,5).Range.Text = '51+30=81'
$t.Cell(18,1).Range.Text = '84-48=36'
$t.Cell(18,2).Range.Text = '78-52=26'
$t.Cell(18,3).Range.Text = '5+2=7'
$t.Cell(18,4).Range.Text = '33-33=0'
$t.Cell(18,5).Range.Text = '82-27=55'
$t.Cell(19,1).Range.Text = '64+28=92'
$t.Cell(19,2).Range.Text = '93-6=87'
$t.Cell(19,3).Range.Text = '2+2=4'
$t.Cell(19,4).Range.Text = '70-28=42'
$t.Cell(19,5).Range.Text = '26+6=32'
$t.Cell(20,1).Range.Text = '44-12=32'
$t.Cell(20,2).Range.Text = '69+13=82'
$t.Cell(20,3).Range.Text = '25+43=68'
$t.Cell(20,4).Range.Text = '26-6=20'
$t.Cell(20,5).Range.Text = '27+36=63'
